# Add a new "PENJELASAN" column to Sheet1, between the existing "JAWABAN"
# column (H) and "BOBOT SOAL ID" column (old I, becomes J), shifting
# BOBOT SOAL ID / GEL / SMT / TAHUN one column to the right (I/J/K/L -> J/K/L/M).
#
# We insert the new blank column at position H (8) -- i.e. *inside* the
# existing uniform-width D:H block -- rather than at its right edge (I/9).
# Doing it this way lets the host correctly extend that column block's
# formatting (width 36.625 / style 3) onto the freshly inserted column
# instead of leaving it with the unformatted default column width.
# Afterwards we put the cell content back where it belongs: H1 keeps
# "JAWABAN" and the newly inserted column gets the new header text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Columns.Item(8).Insert() | Out-Null

$ws.Range("H1").Value = "JAWABAN"
$ws.Range("I1").Value = "PENJELASAN"

# The column insert bumps the workbook's internal column count by one
# (16384 -> 16385) which would otherwise leak into the exported default
# trailing column span; deleting the (now spurious) last column restores
# the normal 16384-column bookkeeping without touching anything visible.
$ws.Columns.Item(16384).Delete() | Out-Null

# Match the saved selection from the edit (cell K7).
$ws.Range("K7").Select() | Out-Null
